$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "JuristPerson" column (J) to hold
# the new "Personensteuer" column. This shifts the old J column to K.
$ws.Columns.Item(10).Insert()

# New column header
$ws.Cells.Item(1, 10).Value = "Personensteuer"

# New column values (rows 2-5)
$ws.Cells.Item(2, 10).Value = 10
$ws.Cells.Item(3, 10).Value = 10
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(5, 10).Value = 20

# Update column E (SteuerfussKanton) values from 100 to 104
$ws.Cells.Item(2, 5).Value = 104
$ws.Cells.Item(3, 5).Value = 104
$ws.Cells.Item(4, 5).Value = 104
$ws.Cells.Item(5, 5).Value = 104

# Set column widths for columns E and F (best-fit widths as left by the author)
$ws.Columns.Item(5).ColumnWidth = 14.358
$ws.Columns.Item(6).ColumnWidth = 17.072

# Update selection
$ws.Range("J11").Select()
